$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 06:52"

# 2. Swap the country labels for rows 73 and 74 (Lituania overtakes Bosnia y Herzegovina
#    in the ranking, so the two countries swap rows) and refresh their statistics.
$ws.Range("A73").Value = "Lituania"
$ws.Range("A74").Value = "Bosnia y Herzegovina"

# Row 73 (now Lituania) - updated figures
$ws.Range("B73").Value = 1091
$ws.Range("C73").Value = 21
$ws.Range("D73").Value = 101
$ws.Range("E73").Value = 961
$ws.Range("F73").Value = 14
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 29

# Row 74 (now Bosnia y Herzegovina) - unchanged figures carried over
$ws.Range("B74").Value = 1083
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 236
$ws.Range("E74").Value = 807
$ws.Range("F74").Value = 4
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 40

# 3. Row 8 - Alemania
$ws.Range("D8").Value = 72600
$ws.Range("E8").Value = 56115

# 4. Row 20 - Austria
$ws.Range("B20").Value = 14234
$ws.Range("C20").Value = 8
$ws.Range("E20").Value = 6217

# 5. Row 35 - Chequia
$ws.Range("B35").Value = 6141
$ws.Range("C35").Value = 30
$ws.Range("E35").Value = 5338

# 6. Row 36 - Pakistan
$ws.Range("B36").Value = 5988
$ws.Range("C36").Value = 151
$ws.Range("D36").Value = 1446
$ws.Range("E36").Value = 4435

# 7. Row 53 - Tailandia
$ws.Range("B53").Value = 2643
$ws.Range("C53").Value = 30
$ws.Range("D53").Value = 1497
$ws.Range("E53").Value = 1103
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 43

# 8. Row 98 - Kirguistan
$ws.Range("B98").Value = 449
$ws.Range("C98").Value = 19
$ws.Range("D98").Value = 78
$ws.Range("E98").Value = 366

# 9. Row 134 - Jamaica
$ws.Range("D134").Value = 21
$ws.Range("E134").Value = 79
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 5
